# Applies commit "#7 add controls: SIZE, MTB, LEV; variable creation and data screening"
# to sheet "T1PA_raw": updates summary statistics for existing variables,
# relabels two rows, and appends four new rows (TLAG, NEG, SIZE, MTB, LEV).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T1PA_raw")

# --- Relabel rows whose variable name changed (case normalisation) ---
$ws.Cells.Item(17, 1).Value = "NW"
$ws.Cells.Item(28, 1).Value = "TONE"

# --- Add new rows for the newly-introduced control variables ---
# Copy the formatting of the last existing label cell (A28) onto the new label cells
$ws.Range("A28").Copy()
$ws.Range("A29:A33").PasteSpecial(-4122)
$ws.Cells.Item(29, 1).Value = "TLAG"
$ws.Cells.Item(30, 1).Value = "NEG"
$ws.Cells.Item(31, 1).Value = "SIZE"
$ws.Cells.Item(32, 1).Value = "MTB"
$ws.Cells.Item(33, 1).Value = "LEV"

# --- Write the summary-statistics values (count, mean, std, min, 25%, 50%, 75%, max) ---
$ws.Cells.Item(2, 2).Value = 202882
$ws.Cells.Item(2, 3).Value = 994.17
$ws.Cells.Item(2, 4).Value = 4316.58
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 42.97
$ws.Cells.Item(2, 7).Value = 141.5
$ws.Cells.Item(2, 8).Value = 499.68
$ws.Cells.Item(2, 9).Value = 167633

$ws.Cells.Item(3, 2).Value = 253295
$ws.Cells.Item(3, 3).Value = 7021.52
$ws.Cells.Item(3, 4).Value = 61050.65
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 112.55
$ws.Cells.Item(3, 7).Value = 475.07
$ws.Cells.Item(3, 8).Value = 1952.41
$ws.Cells.Item(3, 9).Value = 2764661

$ws.Cells.Item(4, 2).Value = 253295
$ws.Cells.Item(4, 3).Value = 1450.06
$ws.Cells.Item(4, 4).Value = 7986.68
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 50.45
$ws.Cells.Item(4, 7).Value = 161.43
$ws.Cells.Item(4, 8).Value = 613.34
$ws.Cells.Item(4, 9).Value = 397609

$ws.Cells.Item(5, 2).Value = 253212
$ws.Cells.Item(5, 3).Value = 818.7
$ws.Cells.Item(5, 4).Value = 11245.56
$ws.Cells.Item(5, 5).Value = -26
$ws.Cells.Item(5, 6).Value = 9.71
$ws.Cells.Item(5, 7).Value = 40.22
$ws.Cells.Item(5, 8).Value = 148.17
$ws.Cells.Item(5, 9).Value = 748548

$ws.Cells.Item(6, 2).Value = 253295
$ws.Cells.Item(6, 3).Value = 97.2
$ws.Cells.Item(6, 4).Value = 405.31
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 10.45
$ws.Cells.Item(6, 7).Value = 25.4
$ws.Cells.Item(6, 8).Value = 61.64
$ws.Cells.Item(6, 9).Value = 29206.44

$ws.Cells.Item(7, 2).Value = 253295
$ws.Cells.Item(7, 3).Value = 674.67
$ws.Cells.Item(7, 4).Value = 10687.73
$ws.Cells.Item(7, 5).Value = -0.09
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 2.86
$ws.Cells.Item(7, 8).Value = 32.85
$ws.Cells.Item(7, 9).Value = 519230

$ws.Cells.Item(8, 2).Value = 253295
$ws.Cells.Item(8, 3).Value = 1197.97
$ws.Cells.Item(8, 4).Value = 9904.44
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0.22
$ws.Cells.Item(8, 7).Value = 36.8
$ws.Cells.Item(8, 8).Value = 334.48
$ws.Cells.Item(8, 9).Value = 616814

$ws.Cells.Item(9, 2).Value = 232652
$ws.Cells.Item(9, 3).Value = 32.34
$ws.Cells.Item(9, 4).Value = 174.13
$ws.Cells.Item(9, 5).Value = -33
$ws.Cells.Item(9, 6).Value = 0.52
$ws.Cells.Item(9, 7).Value = 2.43
$ws.Cells.Item(9, 8).Value = 12.52
$ws.Cells.Item(9, 9).Value = 8166

$ws.Cells.Item(10, 2).Value = 252950
$ws.Cells.Item(10, 3).Value = 45.85
$ws.Cells.Item(10, 4).Value = 368.69
$ws.Cells.Item(10, 5).Value = -41847.9
$ws.Cells.Item(10, 6).Value = -0.4
$ws.Cells.Item(10, 7).Value = 2.29
$ws.Cells.Item(10, 8).Value = 15.28
$ws.Cells.Item(10, 9).Value = 22628

$ws.Cells.Item(11, 2).Value = 180278
$ws.Cells.Item(11, 3).Value = 992.09
$ws.Cells.Item(11, 4).Value = 5868.83
$ws.Cells.Item(11, 5).Value = -9.19
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 22.01
$ws.Cells.Item(11, 8).Value = 238.44
$ws.Cells.Item(11, 9).Value = 312576

$ws.Cells.Item(12, 2).Value = 204118
$ws.Cells.Item(12, 3).Value = 693.41
$ws.Cells.Item(12, 4).Value = 3366.15
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 14.89
$ws.Cells.Item(12, 7).Value = 58.12
$ws.Cells.Item(12, 8).Value = 260.72
$ws.Cells.Item(12, 9).Value = 131246

$ws.Cells.Item(13, 2).Value = 219608
$ws.Cells.Item(13, 3).Value = 702.02
$ws.Cells.Item(13, 4).Value = 3366.38
$ws.Cells.Item(13, 5).Value = -3038.25
$ws.Cells.Item(13, 6).Value = 16.42
$ws.Cells.Item(13, 7).Value = 72.55
$ws.Cells.Item(13, 8).Value = 323.45
$ws.Cells.Item(13, 9).Value = 130377

$ws.Cells.Item(14, 2).Value = 199487
$ws.Cells.Item(14, 3).Value = 170.01
$ws.Cells.Item(14, 4).Value = 1283.17
$ws.Cells.Item(14, 5).Value = -705
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 16.92
$ws.Cells.Item(14, 9).Value = 64381

$ws.Cells.Item(15, 2).Value = 253295
$ws.Cells.Item(15, 3).Value = 32.32
$ws.Cells.Item(15, 4).Value = 778.63
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 7.14
$ws.Cells.Item(15, 7).Value = 16.4
$ws.Cells.Item(15, 8).Value = 31.38
$ws.Cells.Item(15, 9).Value = 118510

$ws.Cells.Item(16, 2).Value = 253295
$ws.Cells.Item(16, 3).Value = 0.02
$ws.Cells.Item(16, 4).Value = 0.28
$ws.Cells.Item(16, 5).Value = -1.58
$ws.Cells.Item(16, 6).Value = -0.11
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0.12
$ws.Cells.Item(16, 9).Value = 18.31

$ws.Cells.Item(17, 2).Value = 253295
$ws.Cells.Item(17, 3).Value = 11363.02
$ws.Cells.Item(17, 4).Value = 10213.56
$ws.Cells.Item(17, 5).Value = 1185
$ws.Cells.Item(17, 6).Value = 4913
$ws.Cells.Item(17, 7).Value = 9010
$ws.Cells.Item(17, 8).Value = 14630
$ws.Cells.Item(17, 9).Value = 722159

$ws.Cells.Item(18, 2).Value = 253295
$ws.Cells.Item(18, 3).Value = 1649.79
$ws.Cells.Item(18, 4).Value = 775.9
$ws.Cells.Item(18, 5).Value = 215
$ws.Cells.Item(18, 6).Value = 1111
$ws.Cells.Item(18, 7).Value = 1552
$ws.Cells.Item(18, 8).Value = 2023
$ws.Cells.Item(18, 9).Value = 13330

$ws.Cells.Item(19, 2).Value = 253295
$ws.Cells.Item(19, 3).Value = 188.2
$ws.Cells.Item(19, 4).Value = 232.94
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 51
$ws.Cells.Item(19, 7).Value = 114
$ws.Cells.Item(19, 8).Value = 234
$ws.Cells.Item(19, 9).Value = 9603

$ws.Cells.Item(20, 2).Value = 253295
$ws.Cells.Item(20, 3).Value = 69.29
$ws.Cells.Item(20, 4).Value = 74.21
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 23
$ws.Cells.Item(20, 7).Value = 49
$ws.Cells.Item(20, 8).Value = 89
$ws.Cells.Item(20, 9).Value = 2828

$ws.Cells.Item(21, 2).Value = 253295
$ws.Cells.Item(21, 3).Value = 156.3
$ws.Cells.Item(21, 4).Value = 161.19
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 53
$ws.Cells.Item(21, 7).Value = 111
$ws.Cells.Item(21, 8).Value = 197
$ws.Cells.Item(21, 9).Value = 4248

$ws.Cells.Item(22, 2).Value = 253295
$ws.Cells.Item(22, 3).Value = 135.88
$ws.Cells.Item(22, 4).Value = 296.98
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 29
$ws.Cells.Item(22, 7).Value = 67
$ws.Cells.Item(22, 8).Value = 140
$ws.Cells.Item(22, 9).Value = 27913

$ws.Cells.Item(23, 2).Value = 253295
$ws.Cells.Item(23, 3).Value = 79.69
$ws.Cells.Item(23, 4).Value = 98.58
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 23
$ws.Cells.Item(23, 7).Value = 54
$ws.Cells.Item(23, 8).Value = 102
$ws.Cells.Item(23, 9).Value = 7103

$ws.Cells.Item(24, 2).Value = 253295
$ws.Cells.Item(24, 3).Value = 31.62
$ws.Cells.Item(24, 4).Value = 35.44
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 10
$ws.Cells.Item(24, 7).Value = 22
$ws.Cells.Item(24, 8).Value = 41
$ws.Cells.Item(24, 9).Value = 1804

$ws.Cells.Item(25, 2).Value = 253295
$ws.Cells.Item(25, 3).Value = 35.46
$ws.Cells.Item(25, 4).Value = 34.25
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 12
$ws.Cells.Item(25, 7).Value = 26
$ws.Cells.Item(25, 8).Value = 47
$ws.Cells.Item(25, 9).Value = 992

$ws.Cells.Item(26, 2).Value = 253295
$ws.Cells.Item(26, 3).Value = 63.18
$ws.Cells.Item(26, 4).Value = 98.63
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 14
$ws.Cells.Item(26, 7).Value = 31
$ws.Cells.Item(26, 8).Value = 64
$ws.Cells.Item(26, 9).Value = 3190

$ws.Cells.Item(27, 2).Value = 253295
$ws.Cells.Item(27, 3).Value = 4.37
$ws.Cells.Item(27, 4).Value = 9.8
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 4
$ws.Cells.Item(27, 9).Value = 400

$ws.Cells.Item(28, 2).Value = 253295
$ws.Cells.Item(28, 3).Value = -0.01
$ws.Cells.Item(28, 4).Value = 0.01
$ws.Cells.Item(28, 5).Value = -0.06
$ws.Cells.Item(28, 6).Value = -0.01
$ws.Cells.Item(28, 7).Value = -0.01
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0.02

$ws.Cells.Item(29, 2).Value = 253295
$ws.Cells.Item(29, 3).Value = 40.29
$ws.Cells.Item(29, 4).Value = 17.48
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 37
$ws.Cells.Item(29, 7).Value = 40
$ws.Cells.Item(29, 8).Value = 44
$ws.Cells.Item(29, 9).Value = 4072

$ws.Cells.Item(30, 2).Value = 253295
$ws.Cells.Item(30, 3).Value = 0.5
$ws.Cells.Item(30, 4).Value = 0.5
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 1
$ws.Cells.Item(30, 9).Value = 1

$ws.Cells.Item(31, 2).Value = 253295
$ws.Cells.Item(31, 3).Value = 5.99
$ws.Cells.Item(31, 4).Value = 1.97
$ws.Cells.Item(31, 5).Value = 2.04
$ws.Cells.Item(31, 6).Value = 4.52
$ws.Cells.Item(31, 7).Value = 5.86
$ws.Cells.Item(31, 8).Value = 7.3
$ws.Cells.Item(31, 9).Value = 11.14

$ws.Cells.Item(32, 2).Value = 253295
$ws.Cells.Item(32, 3).Value = 3.38
$ws.Cells.Item(32, 4).Value = 4.74
$ws.Cells.Item(32, 5).Value = 0.34
$ws.Cells.Item(32, 6).Value = 1.23
$ws.Cells.Item(32, 7).Value = 1.96
$ws.Cells.Item(32, 8).Value = 3.49
$ws.Cells.Item(32, 9).Value = 34.56

$ws.Cells.Item(33, 2).Value = 253295
$ws.Cells.Item(33, 3).Value = 0.2
$ws.Cells.Item(33, 4).Value = 0.19
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 0.02
$ws.Cells.Item(33, 7).Value = 0.15
$ws.Cells.Item(33, 8).Value = 0.32
$ws.Cells.Item(33, 9).Value = 0.77

Write-Host "Updated T1PA_raw summary statistics ($($ws.UsedRange.Rows.Count) rows)"